# Generate Report for Handback
#
# This applies the "handback" update to the localization-status workbook:
#   - Status columns flip from "Ready for handoff" to
#     "Handed back: in sync with en-US" (both on the per-locale sheets and
#     on the Overview roll-up sheet, which shares the same status text).
#   - Each per-locale sheet (zh-cn, de-de) grows two new columns:
#       F = Latest Target File    (same link/display as column A)
#       G = Latest Handback File  (same link/display as column D)
#   - Latest Handback DateTime (column H) is stamped with the handback time,
#     per locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview roll-up sheet: columns B (zh-cn) and C (de-de) mirror the same
# per-locale status text shown on the locale sheets.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-locale sheets
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Handback = "2016-03-17 16:49:11" },
    @{ Name = "de-de"; Handback = "2016-03-17 16:49:17" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    foreach ($row in 2,3) {
        # Pull the existing Source (A) / Handback-file (D) hyperlink info so
        # the new Latest Target File (F) / Latest Handback File (G) columns
        # can reuse the same address + display text.
        $aAddress = $null
        $aDisplay = $null
        $dAddress = $null
        $dDisplay = $null

        $aRef = '$A$' + $row
        $dRef = '$D$' + $row

        foreach ($h in $ws.Hyperlinks) {
            $addr = $h.Range.Address()
            if ($addr -eq $aRef) {
                $aAddress = $h.Address
                $aDisplay = $h.TextToDisplay
            }
            if ($addr -eq $dRef) {
                $dAddress = $h.Address
                $dDisplay = $h.TextToDisplay
            }
        }

        # Status (column C)
        $ws.Range("C" + $row).Value = $newStatus

        # Latest Target File (column F) - mirrors column A
        $fCell = $ws.Range("F" + $row)
        $ws.Hyperlinks.Add($fCell, $aAddress, "", "", $aDisplay)
        $fCell.Font.Name = "Calibri"
        $fCell.Font.Size = 11
        $fCell.Font.Underline = 2
        $fCell.Font.Color = 15570276

        # Latest Handback File (column G) - mirrors column D
        $gCell = $ws.Range("G" + $row)
        $ws.Hyperlinks.Add($gCell, $dAddress, "", "", $dDisplay)
        $gCell.Font.Name = "Calibri"
        $gCell.Font.Size = 11
        $gCell.Font.Underline = 2
        $gCell.Font.Color = 15570276

        # Latest Handback DateTime (column H)
        $ws.Range("H" + $row).Value = $locale.Handback
    }
}

Write-Output "Handback report generated"
